# chore(results): Auto-update draw results on excel 2025-11-30T17:37:37Z
#
# Appends the newest Pick 4 draw as row 75 of the "Results" sheet, mirroring
# the existing rows (Date, Game, Phase, Result, InsertedAt). In the source
# file every column is stored as literal text -- including the date-looking
# and number-looking ones -- so a plain `.Value = "2025-11-30"` would get
# auto-converted by Excel into a date serial / number for the Date and Phase
# columns. To keep them as text (matching every other row, and without
# leaving behind a stray number-format / quote-prefix cell style), we enter
# them as a `="..."` text formula and immediately flatten the cell back down
# to its literal value with Copy + PasteSpecial(xlPasteValues).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow     = 75
$date       = "2025-11-30"
$game       = "Pick 4"
$phase      = "251130"
$result     = "5-9-5-6"
$insertedAt = "2025-11-30T21:37:37.269+04:00"

$xlPasteValues = -4163

$dateCell  = $ws.Cells.Item($newRow, 1)
$phaseCell = $ws.Cells.Item($newRow, 3)

$dateCell.Formula = '="' + $date + '"'
$dateCell.Copy() | Out-Null
$dateCell.PasteSpecial($xlPasteValues) | Out-Null

$phaseCell.Formula = '="' + $phase + '"'
$phaseCell.Copy() | Out-Null
$phaseCell.PasteSpecial($xlPasteValues) | Out-Null

$ws.Cells.Item($newRow, 2).Value = $game        # B75 - plain text already, no conversion risk
$ws.Cells.Item($newRow, 4).Value = $result      # D75 - plain text already, no conversion risk
$ws.Cells.Item($newRow, 5).Value = $insertedAt  # E75 - plain text already, no conversion risk

$excel.CutCopyMode = $false
